# Update COVID-19 row (row 6) data values per updated dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = 13550000
$ws.Range("D6").Value = 584000
$ws.Range("E6").Value = 0.072
$ws.Range("F6").Value = 188
$ws.Range("G6").Value = 44028

# Update the active cell selection on the sheet (was H5, now F6) and scroll position
$ws.Range("F6").Select()
